$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the standalone "Search Function (All Users)" paragraph.
#    Deleting its Range (which includes the trailing paragraph mark)
#    merges it away, shifting the following paragraphs up by one.
# ------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ------------------------------------------------------------------
# 2. The paragraph that used to read "Delete Function (All Users)"
#    (bold + underlined) is now paragraph 2. Replace its run text
#    with "Search Function (All Users)" and drop the direct run
#    formatting (bold/underline) so the text becomes plain while the
#    paragraph mark itself keeps the bold/underline that lives in
#    <w:pPr><w:rPr>.
# ------------------------------------------------------------------
$p = $d.Paragraphs(2)
$r = $p.Range
$r.MoveEnd(1, -1)          # exclude the paragraph mark from the range
$r.Delete()                # remove the old run entirely (and its rPr)
$r.InsertAfter("Search Function (All Users)")   # fresh, unformatted run

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the last paragraph
#    ("Bookmark Thumbnail") to the end of the (new) "Search Function
#    (All Users)" paragraph.
#
#    A zero-length bookmark placed exactly on a (non-final) paragraph
#    boundary can get mis-anchored, so we work around that by
#    temporarily inserting a one-character sentinel, wrapping the
#    bookmark around that sentinel (a normal, non-boundary range),
#    and then deleting the sentinel. The bookmark's tracked range
#    collapses correctly in place, right before the paragraph mark.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$target = $d.Paragraphs(2).Range
$target.MoveEnd(1, -1)
$target.InsertAfter("X")
$sentinelStart = $target.End - 1
$sentinelRange = $d.Range($sentinelStart, $target.End)
$d.Bookmarks.Add("_GoBack", $sentinelRange)
$d.Range($sentinelStart, $sentinelStart + 1).Delete()
